$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '328.91'
Set-TextValue 'E2' '-0.34%'
Set-TextValue 'D3' '43.84'
Set-TextValue 'E3' '-0.78%'
Set-TextValue 'D4' '5.362'
Set-TextValue 'E4' '-2.92%'
Set-TextValue 'D5' '0.08396'
Set-TextValue 'E5' '3.47%'
Set-TextValue 'D6' '4.436'
Set-TextValue 'E6' '0.37%'
Set-TextValue 'D7' '1.944'
Set-TextValue 'E7' '-6.91%'
Set-TextValue 'D8' '0.9759'
Set-TextValue 'E8' '1.28%'
Set-TextValue 'E9' '-6.05%'
Set-TextValue 'D10' '0.1126'
Set-TextValue 'E10' '-1.02%'
Set-TextValue 'D11' '0.1905'
Set-TextValue 'E11' '1.04%'
Set-TextValue 'D12' '0.09680'
Set-TextValue 'E12' '-3.00%'
Set-TextValue 'D13' '0.04614'
Set-TextValue 'E13' '-1.09%'
Set-TextValue 'E14' '0.44%'
Set-TextValue 'D15' '0.001294'
Set-TextValue 'E15' '3.56%'
Set-TextValue 'D16' '0.006112'
Set-TextValue 'E16' '3.99%'
Set-TextValue 'B17' 'LEO'
Set-TextValue 'C17' 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue 'D17' '3.402'
Set-TextValue 'E17' '0.68%'
Set-TextValue 'B18' 'BitpandaEcosystemToken'
Set-TextValue 'C18' 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
Set-TextValue 'D18' '0.3329'
Set-TextValue 'E18' '0.54%'
Set-TextValue 'B19' 'MCDex'
Set-TextValue 'C19' 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
Set-TextValue 'D19' '9.029'
Set-TextValue 'E19' '-10.84%'
Set-TextValue 'B20' 'ProBitToken'
Set-TextValue 'C20' 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
Set-TextValue 'D20' '0.1371'
Set-TextValue 'E20' '-2.13%'
Set-TextValue 'B21' 'ZBToken'
Set-TextValue 'C21' 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
Set-TextValue 'D21' '0.2550'
Set-TextValue 'E21' '2.42%'
Set-TextValue 'B22' 'CoinExToken'
Set-TextValue 'C22' 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
Set-TextValue 'D22' '0.04165'
Set-TextValue 'E22' '1.80%'
Set-TextValue 'B23' 'BitKan'
Set-TextValue 'C23' 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
Set-TextValue 'D23' '0.001295'
Set-TextValue 'E23' '-0.73%'
Set-TextValue 'B24' 'HotbitToken'
Set-TextValue 'C24' 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
Set-TextValue 'D24' '0.004403'
Set-TextValue 'E24' '0.78%'
Set-TextValue 'D25' '0.0001302'
Set-TextValue 'E25' '4.81%'
Set-TextValue 'E26' '-19.65%'
Set-TextValue 'D38' '0.02664'
Set-TextValue 'E38' '-0.65%'
Set-TextValue 'D39' '0.05637'
Set-TextValue 'E39' '-0.38%'
Set-TextValue 'D40' '0.007864'
Set-TextValue 'E40' '3.60%'
Set-TextValue 'D41' '0.1417'
Set-TextValue 'E41' '0.53%'
Set-TextValue 'D42' '0.007359'
Set-TextValue 'E42' '0.42%'
Set-TextValue 'D43' '0.002113'
Set-TextValue 'E43' '7.08%'
Set-TextValue 'D44' '0.007908'
Set-TextValue 'E44' '-4.08%'
Set-TextValue 'D45' '0.3519'
Set-TextValue 'D46' '0.00006905'
Set-TextValue 'E46' '-1.85%'
Set-TextValue 'E47' '0.75%'
Set-TextValue 'D48' '0.003508'
Set-TextValue 'E48' '5.45%'
Set-TextValue 'D49' '0.003534'
Set-TextValue 'E49' '41.07%'
Set-TextValue 'E50' '0.75%'
Set-TextValue 'E51' '0.75%'
